$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and report week dates)
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# Crime statistics grid updates
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = -40
$ws.Range("N14").Value = -82.352941176470
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 62.5
$ws.Range("L15").Value = 62.5
$ws.Range("M15").Value = 18.181818181818
$ws.Range("N15").Value = -60.606060606060
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 94
$ws.Range("K16").Value = -28.723404255319
$ws.Range("L16").Value = 52.272727272727
$ws.Range("M16").Value = -35.576923076923
$ws.Range("N16").Value = -79.510703363914
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -21.951219512195
$ws.Range("I17").Value = 167
$ws.Range("J17").Value = 198
$ws.Range("K17").Value = -15.656565656565
$ws.Range("L17").Value = 1.212121212121
$ws.Range("M17").Value = 60.576923076923
$ws.Range("N17").Value = -53.351955307262
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -7.142857142857
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = -27.659574468085
$ws.Range("L18").Value = 3.030303030303
$ws.Range("M18").Value = 83.783783783783
$ws.Range("N18").Value = -81.471389645776
$ws.Range("C19").Value = 9
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 30.434782608695
$ws.Range("I19").Value = 132
$ws.Range("J19").Value = 128
$ws.Range("K19").Value = 3.125
$ws.Range("L19").Value = 22.222222222222
$ws.Range("M19").Value = 11.864406779661
$ws.Range("N19").Value = -8.333333333333
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = -13.888888888888
$ws.Range("L20").Value = 63.157894736842
$ws.Range("M20").Value = 93.75
$ws.Range("N20").Value = -72.072072072072
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 14.814814814814
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 10.752688172043
$ws.Range("I21").Value = 481
$ws.Range("J21").Value = 563
$ws.Range("K21").Value = -14.564831261101
$ws.Range("L21").Value = 16.183574879227
$ws.Range("M21").Value = 23.017902813299
$ws.Range("N21").Value = -64.554163596168
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -40
$ws.Range("L22").Value = 0
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = 50
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 228.571428571429
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 72
$ws.Range("K23").Value = 23.611111111111
$ws.Range("L23").Value = 15.584415584415
$ws.Range("M23").Value = 45.901639344262
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = -32.051282051282
$ws.Range("I24").Value = 309
$ws.Range("J24").Value = 343
$ws.Range("K24").Value = -9.912536443148
$ws.Range("L24").Value = 7.291666666666
$ws.Range("M24").Value = 20.703125
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 77.777777777777
$ws.Range("F25").Value = 65
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 44.444444444444
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 227
$ws.Range("K25").Value = 3.524229074889
$ws.Range("L25").Value = 22.395833333333
$ws.Range("M25").Value = -27.692307692307
$ws.Range("C26").Value = 1
$ws.Range("L26").Value = 12.5
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = -31.034482758620
$ws.Range("D28").Value = 3
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -28.571428571428
$ws.Range("N28").Value = -75.609756097561
$ws.Range("D29").Value = 3
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -38.461538461538
$ws.Range("N29").Value = -78.378378378378
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 1
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("K30").Value = 0
$ws.Range("K30").NumberFormat = "#,##0.0;""-""#,##0.0"
